$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.478.28"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.884.51"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4717"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2886"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06473"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.23"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07767"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "1.888.39"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.72"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7249"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.182"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "281.40"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").Value = "30.475.79"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "2.134.96"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.267"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.290"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.82"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.067"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.888"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.335"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09645"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.471"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.269"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.139"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04854"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6923"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.817"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "74.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.205"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.962"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4265"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8265"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.10"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.656"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.958"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.20"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "904.02"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05748"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.60%  "
